# Selenium_WebDriver_Advanced_Usage_V2 - "End of day 6" edit
#
# Applies the capitalisation fixes ("Webdriver" -> "WebDriver") and the
# title-slide run split that the author made on day 6, using the
# PowerPoint COM object model (TextRange / Characters) so the XML that
# PowerPoint itself would emit is reproduced as closely as possible.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# Slide 10 ("You have successfully completed ...") - Text Placeholder 2
# "Selenium Webdriver – Advanced Usage"
#   -> two runs: "Selenium WebDriver " + "– Advanced Usage"
# ---------------------------------------------------------------------
$slide10 = $p.Slides.Item(10)
$shp10 = $slide10.Shapes.Item(2)
$tr10 = $shp10.TextFrame.TextRange
$tr10.Characters(1, 19).Text = "Selenium WebDriver "

# ---------------------------------------------------------------------
# Slide 4 (Implicit Wait) - TextBox 1
# standalone "Webdriver" code token -> "WebDriver"
# ---------------------------------------------------------------------
$slide4 = $p.Slides.Item(4)
$shp4 = $slide4.Shapes.Item(3)
$tr4 = $shp4.TextFrame.TextRange
$tr4.Paragraphs(9, 1).Characters(1, 9).Text = "WebDriver"

# ---------------------------------------------------------------------
# Slide 6 (Explicit / WebDriverWait example) - TextBox 4
#   "WebdriverWait wait = new WebdriverWait(driver, 10); "
#     -> "WebDriverWait wait = new WebDriverWait(driver, 10); "
#   "WebdriverWait by default calls ..." -> "WebDriverWait by default calls ..."
# ---------------------------------------------------------------------
$slide6 = $p.Slides.Item(6)
$shp6 = $slide6.Shapes.Item(3)
$tr6 = $shp6.TextFrame.TextRange

$para6_5 = $tr6.Paragraphs(5, 1)
$para6_5.Characters(1, 13).Text = "WebDriverWait"
$para6_5.Characters(26, 13).Text = "WebDriverWait"

$tr6.Paragraphs(12, 1).Characters(1, 13).Text = "WebDriverWait"

# ---------------------------------------------------------------------
# Slide 8 (FluentWait syntax) - TextBox 4
# first "wait" (the type name placeholder) -> "Wait"
# ---------------------------------------------------------------------
$slide8 = $p.Slides.Item(8)
$shp8 = $slide8.Shapes.Item(3)
$tr8 = $shp8.TextFrame.TextRange
$tr8.Paragraphs(11, 1).Characters(1, 4).Text = "Wait"

# ---------------------------------------------------------------------
# Slide 9 (FluentWait example) - TextBox 4
# "FluentWait<Webdriver>" -> "FluentWait<WebDriver>" (both occurrences)
# ---------------------------------------------------------------------
$slide9 = $p.Slides.Item(9)
$shp9 = $slide9.Shapes.Item(3)
$tr9 = $shp9.TextFrame.TextRange
$tr9.Paragraphs(3, 1).Characters(1, 21).Text = "FluentWait<WebDriver>"
$tr9.Paragraphs(4, 1).Characters(1, 22).Text = "FluentWait<WebDriver>("
